$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Qty executed upto date" (C) column values (numeric cells)
$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 41
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 66
$ws.Range("C13").Value = 80
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 80
$ws.Range("C16").Value = 79
$ws.Range("C17").Value = 83

# Update the dependent "Upto date Amount" (G) column values.
# These are stored as text in the workbook (e.g. "1280.00"), so force a
# text number format before assigning so Excel keeps them as strings
# instead of coercing to numeric values.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "1280.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19352.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "5958.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "10880.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "138.00"

# Update the Grand Total rows (G19/H19 and G21/H21)
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "37608.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "37608.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "37608.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "37608.00"
